{"js": "// Locate the existing list item that ends with \"...the wrong ones alone\n// together.\" (the first numbered item describing the river-crossing puzzle)\n// and insert a brand-new list item right after it containing the answer\n// about moving the parrot first.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"the wrong ones alone together.\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the anchor paragraph ending in 'the wrong ones alone together.'\");\n}\n\n// insertParagraph after the anchor inherits its paragraph formatting\n// (ListParagraph style + the same numbered-list numPr), matching the\n// target OOXML exactly.\nanchor.insertParagraph(\n  \"The man must take the parrot across the river first because the cat won\\u2019t eat the bag of seed.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Locate the existing list item that ends with \"...the wrong ones alone\n# together.\" (the first numbered item describing the river-crossing puzzle)\n# and insert a brand-new list item right after it containing the answer\n# about moving the parrot first.\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*the wrong ones alone together.*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the anchor paragraph ending in 'the wrong ones alone together.'\"\n}\n\n# Inserting a paragraph mark after the anchor's range creates a brand new\n# paragraph that inherits the anchor's paragraph formatting (the\n# ListParagraph style + the numbered-list numPr), matching the target\n# OOXML exactly.\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n\n# The new paragraph currently contains only its own paragraph mark;\n# shrink the range by one character so assigning .Text replaces the\n# (empty) body without touching the mark, then set the sentence text.\n$r = $newPara.Range\n$r.End = $r.End - 1\n$r.Text = \"The man must take the parrot across the river first because the cat won\" + [char]0x2019 + \"t eat the bag of seed.\"\n"}
